$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 31 (Leve Item ID 4576)
$ws.Range("H31").Value = 167.66667
$ws.Range("I31").Value = 167.66667
$ws.Range("K31").Value = 503.00001
$ws.Range("M31").Value = -273.00001
# row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 399.66666
$ws.Range("I41").Value = 324.5
$ws.Range("J41").Value = 437.25
$ws.Range("K41").Value = 324.5
$ws.Range("L41").Value = 437.25
$ws.Range("M41").Value = 115.5
$ws.Range("N41").Value = -1317.25
# row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 1146
$ws.Range("J58").Value = 1599.5
$ws.Range("L58").Value = 4798.5
$ws.Range("N58").Value = -5098.5
# row 69 (Leve Item ID 12616)
$ws.Range("H69").Value = 260000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 260000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 780000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -781748
# row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 962.6667
$ws.Range("I70").Value = 944
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 2832
$ws.Range("L70").Value = 3000
$ws.Range("M70").Value = -2562
$ws.Range("N70").Value = -3540
# row 72 (Leve Item ID 12616)
$ws.Range("H72").Value = 260000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 260000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 2340000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -2348736
# row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 962.6667
$ws.Range("I73").Value = 944
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 2832
$ws.Range("L73").Value = 3000
$ws.Range("M73").Value = -1896
$ws.Range("N73").Value = -4872
# row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 2627.8572
$ws.Range("J112").Value = 2637.3
$ws.Range("L112").Value = 7911.900000000001
$ws.Range("N112").Value = -10127.9
# row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 69535.63
$ws.Range("I116").Value = 64489.2
$ws.Range("J116").Value = 120000
$ws.Range("K116").Value = 64489.2
$ws.Range("L116").Value = 120000
$ws.Range("M116").Value = -61047.2
$ws.Range("N116").Value = -126884
# row 121 (Leve Item ID 39731)
$ws.Range("H121").Value = 3000
$ws.Range("J121").Value = 3000
$ws.Range("L121").Value = 9000
$ws.Range("N121").Value = -12494

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 8013.7144
$ws.Range("I2").Value = 8149.25
$ws.Range("J2").Value = 7833
$ws.Range("K2").Value = 8149.25
$ws.Range("L2").Value = 7833
$ws.Range("M2").Value = -8036.25
$ws.Range("N2").Value = -8059
# row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 32309.572
$ws.Range("J63").Value = 14853
$ws.Range("L63").Value = 14853
$ws.Range("N63").Value = -16225
# row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 32309.572
$ws.Range("J66").Value = 14853
$ws.Range("L66").Value = 74265
$ws.Range("N66").Value = -81129
# row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 3919.8
$ws.Range("I88").Value = 2300
$ws.Range("K88").Value = 2300
$ws.Range("M88").Value = -1894
# row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 3919.8
$ws.Range("I91").Value = 2300
$ws.Range("K91").Value = 2300
$ws.Range("M91").Value = -896
# row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 3372.25
$ws.Range("I102").Value = 2999
$ws.Range("J102").Value = 3496.6667
$ws.Range("K102").Value = 2999
$ws.Range("L102").Value = 3496.6667
$ws.Range("M102").Value = -1377
$ws.Range("N102").Value = -6740.6667
# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 8013.7144
$ws.Range("I116").Value = 8149.25
$ws.Range("J116").Value = 7833
$ws.Range("K116").Value = 8149.25
$ws.Range("L116").Value = 7833
$ws.Range("M116").Value = -5855.25
$ws.Range("N116").Value = -12421

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 8013.7144
$ws.Range("I3").Value = 8149.25
$ws.Range("J3").Value = 7833
$ws.Range("K3").Value = 8149.25
$ws.Range("L3").Value = 7833
$ws.Range("M3").Value = -8035.25
$ws.Range("N3").Value = -8061
# row 82 (Leve Item ID 11877)
$ws.Range("H82").Value = 38496.625
$ws.Range("I82").Value = 30994
$ws.Range("J82").Value = 40997.5
$ws.Range("K82").Value = 30994
$ws.Range("L82").Value = 40997.5
$ws.Range("M82").Value = -30611
$ws.Range("N82").Value = -41763.5
# row 85 (Leve Item ID 11877)
$ws.Range("H85").Value = 38496.625
$ws.Range("I85").Value = 30994
$ws.Range("J85").Value = 40997.5
$ws.Range("K85").Value = 30994
$ws.Range("L85").Value = 40997.5
$ws.Range("M85").Value = -29668
$ws.Range("N85").Value = -43649.5
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 4125.625
$ws.Range("I86").Value = 917.5
$ws.Range("K86").Value = 917.5
$ws.Range("M86").Value = 205.5
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 4125.625
$ws.Range("I89").Value = 917.5
$ws.Range("K89").Value = 4587.5
$ws.Range("M89").Value = 1028.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 69 (Leve Item ID 11911)
$ws.Range("H69").Value = 9998
$ws.Range("I69").Value = 9998
$ws.Range("K69").Value = 9998
$ws.Range("M69").Value = -9249
# row 72 (Leve Item ID 11911)
$ws.Range("H72").Value = 9998
$ws.Range("I72").Value = 9998
$ws.Range("K72").Value = 29994
$ws.Range("M72").Value = -26250
# row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1444.125
$ws.Range("I107").Value = 1222
$ws.Range("K107").Value = 1222
$ws.Range("M107").Value = 698

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 25 (Leve Item ID 4709)
$ws.Range("H25").Value = 1000
$ws.Range("J25").Value = 1000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3338
# row 30 (Leve Item ID 4709)
$ws.Range("H30").Value = 1000
$ws.Range("J30").Value = 1000
$ws.Range("L30").Value = 3000
$ws.Range("N30").Value = -3204
# row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# row 58 (Leve Item ID 4703)
$ws.Range("H58").Value = 3000
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9256

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 99 (Leve Item ID 19532)
$ws.Range("H99").Value = 6455.3335
$ws.Range("I99").Value = 6455.3335
$ws.Range("K99").Value = 6455.3335
$ws.Range("M99").Value = -4209.3335

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 16500
$ws.Range("J82").Value = 16500
$ws.Range("L82").Value = 16500
$ws.Range("N82").Value = -17222
# row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 16500
$ws.Range("J85").Value = 16500
$ws.Range("L85").Value = 16500
$ws.Range("N85").Value = -18996
# row 106 (Leve Item ID 18713)
$ws.Range("H106").Value = 49071.25
$ws.Range("J106").Value = 49071.25
$ws.Range("L106").Value = 49071.25
$ws.Range("N106").Value = -51595.25
# row 112 (Leve Item ID 25846)
$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 6612.7144
$ws.Range("I62").Value = 3347.25
$ws.Range("J62").Value = 10966.667
$ws.Range("K62").Value = 3347.25
$ws.Range("L62").Value = 10966.667
$ws.Range("M62").Value = -2723.25
$ws.Range("N62").Value = -12214.667
# row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 6612.7144
$ws.Range("I65").Value = 3347.25
$ws.Range("J65").Value = 10966.667
$ws.Range("K65").Value = 16736.25
$ws.Range("L65").Value = 54833.335
$ws.Range("M65").Value = -13616.25
$ws.Range("N65").Value = -61073.335
# row 74 (Leve Item ID 19022)
$ws.Range("H74").Value = 43743
$ws.Range("J74").Value = 43324.332
$ws.Range("L74").Value = 43324.332
$ws.Range("N74").Value = -45196.332
# row 77 (Leve Item ID 19022)
$ws.Range("H77").Value = 43743
$ws.Range("J77").Value = 43324.332
$ws.Range("L77").Value = 129972.996
$ws.Range("N77").Value = -139332.996
# row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 15000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -19340
